$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = -0.07706499999999999
$ws.Range("G2").Value = -0.1033419023136247
$ws.Range("H2").Value = -0.1033419023136247
$ws.Range("I2").Value = -0.1366867869121763
$ws.Range("J2").Value = -0.1366867869121763
$ws.Range("K2").Value = -1.286
$ws.Range("L2").Value = -0.08264781491002571
$ws.Range("U2").Value = 80.09999999999999
$ws.Range("V2").Value = 1.319604612850082
$ws.Range("W2").Value = -0.02058829642948082
$ws.Range("X2").Value = 0.11873807982155
$ws.Range("Y2").Value = -0.1393263762510309
$ws.Range("Z2").Value = 0.4619097793408733
$ws.Range("AA2").Value = -0.1167696857820146
$ws.Range("AB2").Value = 0.08449761883959862
$ws.Range("AC2").Value = -0.2012673046216132
$ws.Range("AD2").Value = 26.351
$ws.Range("AE2").Value = 2.084232021767317
$ws.Range("AF2").Value = 28.43523202176732
$ws.Range("AG2").Value = -51.66476797823267
$ws.Range("AH2").Value = 0.3190122623433939
$ws.Range("AI2").Value = 0.2188415842210042
$ws.Range("AJ2").Value = -5.718145129396121
$ws.Range("AK2").Value = -1.036711697372378
$ws.Range("AL2").Value = 0.058
$ws.Range("AM2").Value = -1.747
$ws.Range("AN2").Value = -51.466796875
$ws.Range("AO2").Value = -46.03448275862069
$ws.Range("AP2").Value = 100.9077499574857
$ws.Range("AQ2").Value = 1.528334287349743

# --- Row 3 updates ---
$ws.Range("D3").Value = 0.0008699999999999999
$ws.Range("G3").Value = -0.00965034965034965
$ws.Range("H3").Value = -0.00965034965034965
$ws.Range("I3").Value = -0.03981815568285457
$ws.Range("J3").Value = -0.03981815568285457
$ws.Range("K3").Value = -0.949
$ws.Range("L3").Value = -0.06636363636363636
$ws.Range("U3").Value = 21.8
$ws.Range("V3").Value = 1.018691588785047
$ws.Range("W3").Value = -0.03664092664092664
$ws.Range("X3").Value = 0.1525200923228871
$ws.Range("Y3").Value = -0.1891610189638137
$ws.Range("Z3").Value = 0.5429624108524423
$ws.Range("AA3").Value = -0.02161976180526059
$ws.Range("AB3").Value = 0.08477018166860331
$ws.Range("AC3").Value = -0.1063899434738639
$ws.Range("AD3").Value = 25.8
$ws.Range("AE3").Value = 2.076998131324101
$ws.Range("AF3").Value = 27.8769981313241
$ws.Range("AG3").Value = 6.0769981313241
$ws.Range("AH3").Value = 0.5657202992972784
$ws.Range("AI3").Value = 0.5155056511018933
$ws.Range("AJ3").Value = 0.2211667410784678
$ws.Range("AK3").Value = 0.1882764347105287
$ws.Range("AL3").Value = 0.047
$ws.Range("AM3").Value = -0.278
$ws.Range("AN3").Value = 27.04402515723271
$ws.Range("AO3").Value = -23.61702127659575
$ws.Range("AP3").Value = 6.370019005580818
$ws.Range("AQ3").Value = 3.992805755395683

# --- Row 4 new row ---
$ws.Range("A4").Value = "Philippines"
$ws.Range("B4").Value = "Paxys, Inc. (PSE:PAX)"
$ws.Range("C4").Value = "Information Services"
$ws.Range("D4").Value = -0.155
$ws.Range("G4").Value = -1.166666666666667
$ws.Range("H4").Value = -1.166666666666667
$ws.Range("I4").Value = -1.236068871498923
$ws.Range("J4").Value = -1.236068871498923
$ws.Range("K4").Value = -0.337
$ws.Range("L4").Value = -0.2674603174603175
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 58.3
$ws.Range("V4").Value = 1.483460559796438
$ws.Range("W4").Value = -0.004535666218034994
$ws.Range("X4").Value = 0.08495606732021299
$ws.Range("Y4").Value = -0.08949173353824799
$ws.Range("Z4").Value = 0.171446441735713
$ws.Range("AA4").Value = -0.2119196097587686
$ws.Range("AB4").Value = 0.08422505601059393
$ws.Range("AC4").Value = -0.2961446657693626
$ws.Range("AD4").Value = 0.551
$ws.Range("AE4").Value = 0.00723389044321603
$ws.Range("AF4").Value = 0.558233890443216
$ws.Range("AG4").Value = -57.74176610955678
$ws.Range("AH4").Value = 0.01400548483853078
$ws.Range("AI4").Value = 0.007358909663642242
$ws.Range("AJ4").Value = 3.131032340749305
$ws.Range("AK4").Value = -3.288586225120573
$ws.Range("AL4").Value = 0.011
$ws.Range("AM4").Value = -1.469
$ws.Range("AN4").Value = -0.3758526603001365
$ws.Range("AO4").Value = -141.8181818181818
$ws.Range("AP4").Value = 39.38728929710558
$ws.Range("AQ4").Value = 1.061946902654867
